# GSBC (GSBD) SEC filing links: rewrite each txtFileLink (column P) URL so it
# points at the accession's own EDGAR directory (.../data/CIK/ACCESSIONNODASH/ACCESSION.txt)
# instead of directly at .../data/CIK/ACCESSION.txt, and widen column P to fit
# the longer links.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# txtFileLink values live in column P (16), data rows 2-37.
$firstRow = 2
$lastRow = 37
$col = 16

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $old = $cell.Value2

    if ($old -match '^(?<prefix>.*/edgar/data/\d+/)(?<accession>\d{10}-\d{2}-\d{6})\.txt$') {
        $prefix = $matches['prefix']
        $accession = $matches['accession']
        $accessionNoDash = $accession.Replace('-', '')
        $newUrl = $prefix + $accessionNoDash + '/' + $accession + '.txt'
        $cell.Value = $newUrl
    }
}

# Widen column P (16) to fit the now-longer links. Excel's ColumnWidth setter
# snaps to whole-pixel increments, so use the value that lands closest to the
# target stored width (~96.71 chars).
$ws.Columns.Item($col).ColumnWidth = 95.8333333333333
